{"js": "// Color the \"Gesti\u00f3n de promociones\" and \"Gesti\u00f3n de pedidos a proveedores\"\n// heading paragraphs (bold title run, trailing \":\" run, and the paragraph\n// mark itself) in red (EE0000), matching the GestionCompras base\n// implementation commit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Gesti\u00f3n de promociones:\",\n  \"Gesti\u00f3n de pedidos a proveedores:\"\n];\n\nfor (const paragraph of paragraphs.items) {\n  const text = (paragraph.text || \"\").trim();\n  if (targets.indexOf(text) !== -1) {\n    // Setting color on the paragraph's font applies it to the paragraph\n    // mark run properties (w:pPr/w:rPr) as well as every run contained in\n    // the paragraph (the bold \"Gesti\u00f3n de ...\" run and the \":\" run).\n    paragraph.font.color = \"#EE0000\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Color the \"Gesti\u00f3n de promociones\" and \"Gesti\u00f3n de pedidos a proveedores\"\n# heading paragraphs (bold title run, trailing \":\" run, and the paragraph\n# mark itself) in red (EE0000), matching the GestionCompras base\n# implementation commit.\n\n$d = $word.ActiveDocument\n\n$targets = @(\"Gesti\u00f3n de promociones:\", \"Gesti\u00f3n de pedidos a proveedores:\")\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($targets -contains $t) {\n        # Setting Font.Color on the paragraph's Range colors every run in\n        # the paragraph (the bold \"Gesti\u00f3n de ...\" run and the \":\" run) and\n        # also the paragraph mark run properties (w:pPr/w:rPr).\n        $p.Range.Font.Color = \"EE0000\"\n    }\n}\n"}
